$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "64.767.47"
$ws.Range("E2").Value = "  -1.35%  "
Set-TextValue "D3" "3.340.34"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "555.79"
$ws.Range("E5").Value = "  -1.09%  "
Set-TextValue "D6" "172.81"
$ws.Range("E6").Value = "  -1.36%  "
Set-TextValue "D7" "0.625"
$ws.Range("E7").Value = "  -0.20%  "
Set-TextValue "D8" "3.330.64"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("E9").Value = "  -0.05%  "
Set-TextValue "D10" "0.170"
$ws.Range("E10").Value = "  -1.30%  "
Set-TextValue "D11" "0.626"
$ws.Range("E11").Value = "  -0.76%  "
Set-TextValue "D12" "52.69"
$ws.Range("E12").Value = "  -4.24%  "
Set-TextValue "D13" "0.0000273"
$ws.Range("E13").Value = "  -2.82%  "
Set-TextValue "D14" "9.09"
$ws.Range("E14").Value = "  -0.61%  "
Set-TextValue "D15" "3.914.50"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.388.23"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D17" "0.119"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D18" "18.06"
$ws.Range("E18").Value = "  -1.52%  "
Set-TextValue "D19" "64.798.30"
$ws.Range("E19").Value = "  -1.25%  "
Set-TextValue "D20" "11.69"
$ws.Range("E20").Value = "  -1.80%  "
Set-TextValue "D21" "0.986"
$ws.Range("E21").Value = "  -0.92%  "
Set-TextValue "D22" "478.85"
$ws.Range("E22").Value = "  +1.86%  "
Set-TextValue "D23" "4.91"
$ws.Range("E23").Value = "  -3.04%  "
Set-TextValue "D24" "89.20"
$ws.Range("E24").Value = "  +3.20%  "
Set-TextValue "D25" "14.19"
$ws.Range("E25").Value = "  +4.27%  "
Set-TextValue "D26" "4.05"
$ws.Range("E26").Value = "  -2.13%  "
Set-TextValue "D27" "2.87"
$ws.Range("E27").Value = "  -0.61%  "
Set-TextValue "D28" "10.49"
$ws.Range("E28").Value = "  -3.77%  "
Set-TextValue "D29" "8.60"
$ws.Range("E29").Value = "  -3.24%  "
Set-TextValue "D30" "30.89"
$ws.Range("E30").Value = "  +0.19%  "
Set-TextValue "D31" "6.47"
$ws.Range("E31").Value = "  -3.31%  "
Set-TextValue "D32" "62.70"
$ws.Range("E32").Value = "  +1.59%  "
Set-TextValue "D33" "11.33"
$ws.Range("E33").Value = "  -1.90%  "
Set-TextValue "D34" "570.50"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E37").Value = "  +3.44%  "
$ws.Range("E38").Value = "  +0.02%  "
Set-TextValue "D39" "35.20"
$ws.Range("E39").Value = "  -1.88%  "
Set-TextValue "D40" "0.370"
$ws.Range("E40").Value = "  -1.24%  "
Set-TextValue "D41" "0.0₃0728"
$ws.Range("E41").Value = "  -3.91%  "
Set-TextValue "D42" "3.092.91"
$ws.Range("E42").Value = "  -0.27%  "
Set-TextValue "D43" "2.77"
$ws.Range("E43").Value = "  -2.94%  "
Set-TextValue "D44" "0.0410"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("E45").Value = "  -1.43%  "
Set-TextValue "D46" "3.14"
$ws.Range("E46").Value = "  -2.83%  "
Set-TextValue "D47" "2.40"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("E48").Value = "  +0.21%  "
Set-TextValue "D49" "140.52"
$ws.Range("E49").Value = "  +3.25%  "
Set-TextValue "D50" "2.56"
$ws.Range("E50").Value = "  -1.05%  "
Set-TextValue "D51" "8.32"
$ws.Range("E51").Value = "  -0.31%  "
